# Fruta / hortaliza, semanal
# Insert a new weekly price row for "Ajo" (Macroferia Regional de Talca)
# right after the existing row 39, pushing every subsequent row down by
# one (the former row 160 becomes row 161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..160 down to 41..161, leaving a blank row 40 to fill in.
$ws.Rows.Item(40).Insert()

$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44459
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = 100112003
$ws.Range("G40").Value = "Ajo"
$ws.Range("H40").Value = "Chino"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = 15000
$ws.Range("N40").Value = "`$/malla 10 kilos"
$ws.Range("O40").Value = "Región del Maule"
$ws.Range("P40").Value = 1500
$ws.Range("Q40").Value = 10
$ws.Range("R40").Value = "Hortaliza"
